# Remove the "Non myelinating Schwann cells" row from the marker table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 holds: Brain | Non myelinating Schwann cells | SOX10,GAP43,NCAM1,NGFR | | Non myelinating Schwann
$ws.Rows.Item(16).Delete()

# Update the cursor/selection to match the saved state after the edit.
$ws.Range("J14").Select()
